# #5: property aircraft done
# Update the "property_category" column on the 建物 (building) and
# 汽車 (car) sheets so each row reflects its own category instead of
# the previously-copied "land" value.

$wb = $excel.ActiveWorkbook

# 建物 (building) sheet: column I is property_category, rows 2-5.
$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I5").Value = "building"

# 汽車 (car) sheet: column H is property_category, row 2.
$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
